# Update cryptos list values (Price and Volume(1h) columns) with latest scrape data.
# Price strings that are unambiguous (already contain 2+ "." separators, e.g. "41.682.98")
# are plain text in Excel automatically. Price strings that look like an ordinary decimal
# number (e.g. "92.57") need to be forced to Text so Excel keeps them as strings, matching
# the original inline-string cell contents; the cell style is restored right after so no
# visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.682.98"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.477.32"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +8.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "2.860.72"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").Value = "2.469.70"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.794"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "41.626.42"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "1.999.10"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.83%  "
$ws.Range("D48").Value = "2.718.19"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.98%  "
